$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name) from "Through 2022-10-25" to "Through 2022-10-26"
$ws.Name = "Through 2022-10-26"

# Update header text in I1 from "2022 (through 10-25)" to "2022 (through 10-26)"
$ws.Range("I1").Value = "2022 (through 10-26)"

# Update data values: I11 94 -> 95, I14 1371 -> 1372
$ws.Range("I11").Value = 95
$ws.Range("I14").Value = 1372
